$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.770.46"
$ws.Range("E2").Value = "  -3.76%  "
$ws.Range("D3").Value = "2.537.90"
$ws.Range("E3").Value = "  -4.66%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'513.28"
$ws.Range("E5").Value = "  -2.30%  "
$ws.Range("D6").Value = "'140.29"
$ws.Range("E6").Value = "  -2.84%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.555"
$ws.Range("E8").Value = "  -2.52%  "
$ws.Range("D9").Value = "'6.51"
$ws.Range("E9").Value = "  -6.98%  "
$ws.Range("D10").Value = "'0.0992"
$ws.Range("E10").Value = "  -3.69%  "
$ws.Range("E11").Value = "  -3.77%  "
$ws.Range("E12").Value = "  -0.32%  "
$ws.Range("D13").Value = "2.987.28"
$ws.Range("D14").Value = "56.796.25"
$ws.Range("E14").Value = "  -3.72%  "
$ws.Range("D15").Value = "'19.99"
$ws.Range("E15").Value = "  -5.06%  "
$ws.Range("E16").Value = "  -3.29%  "
$ws.Range("D17").Value = "2.510.86"
$ws.Range("E17").Value = "  -5.21%  "
$ws.Range("D18").Value = "'331.58"
$ws.Range("E18").Value = "  -2.17%  "
$ws.Range("D19").Value = "'4.27"
$ws.Range("E19").Value = "  -2.95%  "
$ws.Range("E20").Value = "  -2.91%  "
$ws.Range("E21").Value = "  -4.21%  "
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").Value = "'64.42"
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("E24").Value = "  -0.43%  "
$ws.Range("D25").Value = "'0.998"
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").Value = "'0.399"
$ws.Range("E26").Value = "  -4.49%  "
$ws.Range("D27").Value = "2.653.31"
$ws.Range("E27").Value = "  -4.63%  "
$ws.Range("D28").Value = "'6.85"
$ws.Range("E28").Value = "  -3.59%  "
$ws.Range("E29").Value = "  -6.76%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").Value = "'6.25"
$ws.Range("E31").Value = "  -6.68%  "
$ws.Range("E32").Value = "  -3.04%  "
$ws.Range("D33").Value = "'18.48"
$ws.Range("E33").Value = "  -2.03%  "
$ws.Range("D34").Value = "'147.99"
$ws.Range("E34").Value = "  -1.77%  "
$ws.Range("D35").Value = "'3.98"
$ws.Range("E35").Value = "  -4.27%  "
$ws.Range("E36").Value = "  -5.19%  "
$ws.Range("D37").Value = "'0.841"
$ws.Range("E37").Value = "  -5.93%  "
$ws.Range("D38").Value = "'35.57"
$ws.Range("E38").Value = "  -3.58%  "
$ws.Range("D39").Value = "'0.818"
$ws.Range("E39").Value = "  -6.44%  "
$ws.Range("E40").Value = "  -2.76%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("E42").Value = "  -3.29%  "
$ws.Range("E43").Value = "  -1.93%  "
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("E45").Value = "  -6.30%  "
$ws.Range("D46").Value = "'260.32"
$ws.Range("E46").Value = "  -5.74%  "
$ws.Range("D47").Value = "'0.0517"
$ws.Range("E47").Value = "  -2.65%  "
$ws.Range("D48").Value = "'18.49"
$ws.Range("E48").Value = "  -7.06%  "
$ws.Range("D49").Value = "1.964.40"
$ws.Range("E49").Value = "  -4.27%  "
$ws.Range("E50").Value = "  -4.02%  "
$ws.Range("D51").Value = "'4.50"
$ws.Range("E51").Value = "  -4.33%  "
